$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 43 used to (incorrectly) duplicate "Dragon de Komodo" (row 8). The
# insertion script now points it at the newly-added species, "Gecko léopard",
# and its status_id (C43) moves from 5 to 7.
$ws.Range("A43").Value = "Gecko léopard"
$ws.Range("C43").Value = 7

# Row 56 was a duplicate "Tapir terrestre" entry (row 22 already has it).
# Deleting it shifts every following row up by one.
$ws.Rows(56).Delete()

# After the shift, several status_id values (column C) need to be corrected
# to reflect the removal of that row.
$ws.Range("C51").Value = 9
$ws.Range("C52").Value = 9
$ws.Range("C53").Value = 9
$ws.Range("C54").Value = 9
$ws.Range("C55").Value = 9
$ws.Range("C56").Value = 9
$ws.Range("C57").Value = 9
$ws.Range("C58").Value = 9
$ws.Range("C59").Value = 9
$ws.Range("C60").Value = 9
$ws.Range("C62").Value = 9
$ws.Range("C63").Value = 9
$ws.Range("C64").Value = 9
$ws.Range("C65").Value = 9
$ws.Range("C66").Value = 9
$ws.Range("C67").Value = 9
$ws.Range("C68").Value = 9
$ws.Range("C71").Value = 8
$ws.Range("C73").Value = 9
$ws.Range("C79").Value = 8

# Restore the working-selection state (row 56, the new current row) and the
# scroll position shown while editing.
$ws.Rows(56).Select() | Out-Null
$win = $excel.ActiveWindow
$win.ScrollRow = 34
$win.ScrollColumn = 1 | Out-Null
